$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two new rows at position 14. Everything from the old row 14
#    downward (old rows 14,15,16,18,19,21,22,23,24,25,26,27,29,33) shifts
#    down by exactly two, landing on 16,17,18,20,21,23,24,25,26,27,28,29,31,35
#    -- matching the target layout.
# ---------------------------------------------------------------------------
$ws.Rows("14:15").Insert()

# ---------------------------------------------------------------------------
# 2. New row 15 becomes what row 13 used to be (TERMS OF SALE: / TERMS OF
#    PAYMENT: labels), so copy row 13's formatting down, then overwrite the
#    text of row 13 itself with the new PACKING:/SHIPPING SCHEDULE: labels.
# ---------------------------------------------------------------------------
$ws.Range("A13:N13").Copy()
$ws.Range("A15:N15").PasteSpecial(-4122)
$ws.Range("C15").Value = "TERMS OF SALE:"
$ws.Range("G15").Value = "TERMS OF PAYMENT:"

# ---------------------------------------------------------------------------
# 3. New row 14 becomes what row 12 used to be (the Quote.Packing__c /
#    Quote.Shipping_Schedule__c values), so copy row 12's formatting + values
#    down (text is identical, so a plain paste-all works).
# ---------------------------------------------------------------------------
$ws.Range("A12:N12").Copy()
$ws.Range("A14:N14").PasteSpecial(-4122)
$ws.Range("C14").Value = "{{Quote.Packing__c}}"
$ws.Range("G14").Value = "{{Quote.Shipping_Schedule__c}}"
$ws.Rows("14").RowHeight = 45

# ---------------------------------------------------------------------------
# 4. Row 13 keeps its old formatting, only the labels change.
# ---------------------------------------------------------------------------
$ws.Range("C13").Value = "PACKING:"
$ws.Range("G13").Value = "SHIPPING SCHEDULE:"

# ---------------------------------------------------------------------------
# 5. Rows 11 and 12 become blank (all cells styled like the rest of the
#    merged block, style index used by e.g. D17/C17). Clear values first,
#    then copy that blank formatting across.
# ---------------------------------------------------------------------------
$ws.Range("C11:N12").ClearContents()
$ws.Range("C17:N17").Copy()
$ws.Range("C11:N11").PasteSpecial(-4122)
$ws.Range("C12:N12").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 6. Row 10: D10:F10 and H10:N10 take on the blank-box style (same as
#    C11:N12 above); G10 becomes a checkbox-style cell (same formatting as
#    C16/G16) with the new FREIGHT/PREPAID/COLLECT text.
# ---------------------------------------------------------------------------
$ws.Range("C17:N17").Copy()
$ws.Range("D10:F10").PasteSpecial(-4122)
$ws.Range("H10:N10").PasteSpecial(-4122)

$ws.Range("G16").Copy()
$ws.Range("G10").PasteSpecial(-4122)
$ws.Range("G10").Value = "FREIGHT: {{#if Quote.Incoterms__c 'contains' 'FREIGHT'}}" + [char]9745 + "{{else}}" + [char]9744 + "{{/if}}`nPREPAID: {{#if Quote.Incoterms__c 'contains' 'PREPAID'}}" + [char]9745 + "{{else}}" + [char]9744 + "{{/if}}`nCOLLECT: {{#if Quote.Incoterms__c 'contains' 'COLLECT'}}" + [char]9745 + "{{else}}" + [char]9744 + "{{/if}}"

# ---------------------------------------------------------------------------
# 7. Re-merge the A10 block cells that are now contiguous blank ranges.
# ---------------------------------------------------------------------------
$ws.Range("C10:F12").Merge()
$ws.Range("G10:N12").Merge()
$ws.Range("C14:F14").Merge()
$ws.Range("G14:N14").Merge()
$ws.Range("C15:F15").Merge()
$ws.Range("G15:N15").Merge()
